$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed cryptocurrency Price (D) / Volume(1h) (E) figures.
#
# Prices are stored as literal text in the source sheet (e.g. '1.007',
# '0.00001056', '27.768.38' -- note some use '.' as a thousands separator
# so they are not even valid numbers). Assigning such strings straight to
# .Value lets Excel's input parser auto-coerce the numeric-looking ones to
# a Double, which would corrupt values like '6.270' (-> 6.27) or
# '0.00001056' (-> 1.056E-05 in General format). To force literal text entry
# (with no side-effect on the cell's style/number format) we write the text
# as a quoted-string formula and immediately collapse it to a static value
# via Copy + PasteSpecial(xlPasteValues).
$rowData = @{
    2 = @("27.768.38", "  -0.31%  ")
    3 = @("1.766.15", "  -2.44%  ")
    4 = @("1.007", "  +0.65%  ")
    5 = @("338.18", "  +0.37%  ")
    6 = @($null, "  +0.51%  ")
    7 = @("0.3776", "  -3.77%  ")
    8 = @("0.3372", "  -3.24%  ")
    9 = @("45.57", "  -5.25%  ")
    10 = @("1.131", "  -5.78%  ")
    11 = @("0.07286", "  -3.78%  ")
    12 = @("23.12", "  +4.60%  ")
    13 = @($null, "  +0.47%  ")
    14 = @("6.270", "  -3.73%  ")
    15 = @("7.270", "  +1.17%  ")
    16 = @("1.768.79", "  -2.30%  ")
    17 = @("0.00001056", "  -4.51%  ")
    18 = @("0.06609", "  -1.22%  ")
    19 = @("81.12", "  -4.79%  ")
    20 = @("1.002", "  +0.63%  ")
    21 = @("17.21", "  -3.66%  ")
    22 = @("6.352", "  -3.22%  ")
    23 = @("27.822.36", $null)
    24 = @("11.82", "  -8.09%  ")
    25 = @("2.386", "  -0.81%  ")
    26 = @("1.516", "  +2.80%  ")
    27 = @("20.12", "  -5.45%  ")
    28 = @("151.69", "  -1.96%  ")
    29 = @($null, "  -7.68%  ")
    30 = @("1.968.34", "  -2.69%  ")
    31 = @("133.28", "  -1.77%  ")
    32 = @("4.039", "  +0.12%  ")
    33 = @("5.935", "  -2.97%  ")
    34 = @("0.08778", "  -0.65%  ")
    35 = @("12.43", "  -6.30%  ")
    36 = @("0.02366", "  -2.56%  ")
    37 = @("0.6718", "  -2.71%  ")
    38 = @("0.06278", "  -4.26%  ")
    39 = @("5.209", "  -5.80%  ")
    40 = @("0.2123", "  -4.61%  ")
    43 = @("8.073", "  -5.90%  ")
    44 = @("1.001", "  +0.58%  ")
    45 = @("13.87", "  -5.99%  ")
    46 = @("0.6136", "  -6.41%  ")
    47 = @("3.847", "  -0.45%  ")
    48 = @("131.91", "  -0.38%  ")
    49 = @("2.032", "  -5.76%  ")
    50 = @("0.07284", "  +1.14%  ")
    51 = @("1.189", "  +2.39%  ")
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    $priceText = $vals[0]
    $volumeText = $vals[1]
    if ($priceText -ne $null) {
        $priceCell = $ws.Cells.Item($r, 4)
        $priceCell.Formula = "=""" + $priceText + """"
        $priceCell.Copy()
        $priceCell.PasteSpecial(-4163)
    }
    if ($volumeText -ne $null) {
        $ws.Cells.Item($r, 5).Value = $volumeText
    }
}

# Rows 41 and 42 swapped identity: TrustWalletToken now ranks 41st (was 42nd),
# WEMIXTOKEN drops to 42nd (was 41st), both with freshly refreshed figures.
$ws.Cells.Item(41, 2).Value = "TrustWalletToken"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$price41 = $ws.Cells.Item(41, 4)
$price41.Formula = "=""1.224"""
$price41.Copy()
$price41.PasteSpecial(-4163)
$ws.Cells.Item(41, 5).Value = "  -3.16%  "

$ws.Cells.Item(42, 2).Value = "WEMIXTOKEN"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$price42 = $ws.Cells.Item(42, 4)
$price42.Formula = "=""1.476"""
$price42.Copy()
$price42.PasteSpecial(-4163)
$ws.Cells.Item(42, 5).Value = "  -8.36%  "

$excel.CutCopyMode = 0
